$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "Lucia"
$ws.Range("G3").Value = "Rodolfo"
$ws.Range("K3").Value = "Daiana"
$ws.Range("C4").Value = "Rodolfo"
$ws.Range("D4").Value = "Ediane"
$ws.Range("F4").Value = "Patricia Dias"
$ws.Range("G4").Value = "Robson"
$ws.Range("J4").Value = "Lucia"
$ws.Range("K4").Value = "Vanda"
$ws.Range("D5").Value = "Daniel"
$ws.Range("C6").Value = "Isabele"
$ws.Range("G6").Value = "Edith"
$ws.Range("J6").Value = "Vinicius"
$ws.Range("K6").Value = "Antonio"
$ws.Range("K7").Value = "Eliane"
$ws.Range("M7").Value = "Carlos Eduardo"
$ws.Range("D9").Value = "Lucia"
$ws.Range("F9").Value = "Graca"
$ws.Range("G9").Value = "Aline Lima"
$ws.Range("J9").Value = "Icaro"
$ws.Range("K9").Value = "Vinicius"
$ws.Range("M9").Value = "Marcio"
$ws.Range("D10").Value = "Valquiria"
$ws.Range("G10").Value = "Patricia Rodrigues"
$ws.Range("J10").Value = "Icaro"
$ws.Range("C11").Value = "Rodolfo"
$ws.Range("D11").Value = "Valquiria"
$ws.Range("G11").Value = "Lurdes"
$ws.Range("D12").Value = "Vinicius"
$ws.Range("C13").Value = "Keila"
$ws.Range("D13").Value = "Eliane"
$ws.Range("G13").Value = "Lucia"
$ws.Range("J13").Value = "Beth"
$ws.Range("K13").Value = "Valquiria"
$ws.Range("L13").Value = "Rodolfo"
$ws.Range("M14").Value = "Douglas Oliveira"
$ws.Range("C16").Value = "Alana"
$ws.Range("D16").Value = "Helaine Camilo"
$ws.Range("F16").Value = "Lurdes"
$ws.Range("J16").Value = "Keila"
$ws.Range("K16").Value = "Eliane"
$ws.Range("L16").Value = "Vanda"
$ws.Range("M16").Value = "EMPTY"
$ws.Range("D17").Value = "Eliane"
$ws.Range("G17").Value = "Alana"
$ws.Range("J17").Value = "Rodolfo"
$ws.Range("K17").Value = "EMPTY"
$ws.Range("C18").Value = "Ediane"
$ws.Range("D18").Value = "Lindoia"
$ws.Range("F18").Value = "Edith"
$ws.Range("G18").Value = "Patricia Dias"
$ws.Range("J18").Value = "Lucia"
$ws.Range("K18").Value = "Vanda"
$ws.Range("D19").Value = "Alana"
$ws.Range("C20").Value = "Edith"
$ws.Range("D20").Value = "Isabele"
$ws.Range("G20").Value = "Patricia Rodrigues"
$ws.Range("J20").Value = "Aline Lima"
$ws.Range("K20").Value = "Daiana"
$ws.Range("L20").Value = "Icaro"
$ws.Range("M20").Value = "Clayton"
$ws.Range("M21").Value = "Carlos Eduardo"
$ws.Range("D23").Value = "Lucia"
$ws.Range("F23").Value = "Aline Lima"
$ws.Range("G23").Value = "Graca"
$ws.Range("J23").Value = "Patricia Dias"
$ws.Range("K23").Value = "Valquiria"
$ws.Range("L23").Value = "Antonio"
$ws.Range("M23").Value = "Icaro"
$ws.Range("D24").Value = "Cida"
$ws.Range("J24").Value = "Keila"
$ws.Range("K24").Value = "Jessica Silva"
$ws.Range("C25").Value = "Rodolfo"
$ws.Range("D25").Value = "Aline Lima"
$ws.Range("F25").Value = "Lurdes"
$ws.Range("G25").Value = "Valquiria"
$ws.Range("J25").Value = "Antonio"
$ws.Range("D26").Value = "Lurdes"
$ws.Range("C27").Value = "Lurdes"
$ws.Range("D27").Value = "Keila"
$ws.Range("G27").Value = "Lucia"
$ws.Range("J27").Value = "Valquiria"
$ws.Range("K27").Value = "Beth"
$ws.Range("L27").Value = "Rodolfo"
$ws.Range("M27").Value = "Vinicius"
$ws.Range("M28").Value = "Geronimo"
$ws.Range("C30").Value = "Helaine Camilo"
$ws.Range("F30").Value = "Alana"
$ws.Range("J30").Value = "Keila"
$ws.Range("K30").Value = "Eliane"
$ws.Range("L30").Value = "Dario"
$ws.Range("M30").Value = "Amintas"
$ws.Range("D31").Value = "Vanda"
$ws.Range("G31").Value = "Alana"
$ws.Range("J31").Value = "Daiana"
$ws.Range("D32").Value = "Lucia"
$ws.Range("F32").Value = "Patricia Dias"
$ws.Range("G32").Value = "Robson"
$ws.Range("K32").Value = "EMPTY"
$ws.Range("D33").Value = "Cida"

# Update generation date/time (Data de geracao)
$ws.Range("C35").Value = 43818.62945150064
